$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The site was re-scraped at 2025-11-26 12:38:38 JST; refresh the listing rows
# (existing hyperlinks are removed first and re-created after the rewrite so
# relationship ids line up with the new row order).
$ws.Hyperlinks.Delete()
$ws.Range("A2:H15").ClearContents()

# Row 2
$ws.Cells.Item(2, 1).Value = "2025-11-26 12:38:38"
$ws.Cells.Item(2, 2).Value = "Google AI studio が生成したウェブアプリの調整【AI文章での提案は却下します】"
$ws.Cells.Item(2, 3).Value = "システム開発"
$ws.Cells.Item(2, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(2, 5).Value = "期限情報なし"
$ws.Cells.Item(2, 6).Value = "https://www.lancers.jp/work/detail/5441470"
$ws.Cells.Item(2, 7).Value = 333
$ws.Cells.Item(2, 8).Value = "🔥AI,Ai ◇アプリ"

# Row 3
$ws.Cells.Item(3, 1).Value = "2025-11-26 12:38:38"
$ws.Cells.Item(3, 2).Value = "【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪"
$ws.Cells.Item(3, 3).Value = "システム開発"
$ws.Cells.Item(3, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(3, 5).Value = "期限情報なし"
$ws.Cells.Item(3, 6).Value = "https://www.lancers.jp/work/detail/5217096"
$ws.Cells.Item(3, 7).Value = 243
$ws.Cells.Item(3, 8).Value = "🔥API ◆ツール"

# Row 4
$ws.Cells.Item(4, 1).Value = "2025-11-26 12:38:38"
$ws.Cells.Item(4, 2).Value = "【高単価業務自動化】行政書士向けシステム開発依頼"
$ws.Cells.Item(4, 3).Value = "システム開発"
$ws.Cells.Item(4, 4).Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Cells.Item(4, 5).Value = "期限情報なし"
$ws.Cells.Item(4, 6).Value = "https://www.lancers.jp/work/detail/5441252"
$ws.Cells.Item(4, 7).Value = 210
$ws.Cells.Item(4, 8).Value = "◆開発,システム開発"

# Row 5
$ws.Cells.Item(5, 1).Value = "2025-11-26 12:38:38"
$ws.Cells.Item(5, 2).Value = "初回 FastAPIバックエンドの軽微な修正・調整対応エンジニア募集"
$ws.Cells.Item(5, 3).Value = "システム開発"
$ws.Cells.Item(5, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(5, 5).Value = "期限情報なし"
$ws.Cells.Item(5, 6).Value = "https://www.lancers.jp/work/detail/5441207"
$ws.Cells.Item(5, 7).Value = 183
$ws.Cells.Item(5, 8).Value = "🔥API"

# Row 6
$ws.Cells.Item(6, 1).Value = "2025-11-26 12:38:38"
$ws.Cells.Item(6, 2).Value = "【急募】WEBサイト・アプリ開発に強いコーダーを探しています!"
$ws.Cells.Item(6, 3).Value = "システム開発"
$ws.Cells.Item(6, 4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(6, 5).Value = "期限情報なし"
$ws.Cells.Item(6, 6).Value = "https://www.lancers.jp/work/detail/5441442"
$ws.Cells.Item(6, 7).Value = 135
$ws.Cells.Item(6, 8).Value = "◆開発 ◇アプリ"

# Row 7
$ws.Cells.Item(7, 1).Value = "2025-11-26 12:38:38"
$ws.Cells.Item(7, 2).Value = "【急募】GitHub管理のBootstrapサイト移行作業及びWordPress移行"
$ws.Cells.Item(7, 3).Value = "システム開発"
$ws.Cells.Item(7, 4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(7, 5).Value = "期限情報なし"
$ws.Cells.Item(7, 6).Value = "https://www.lancers.jp/work/detail/5441612"
$ws.Cells.Item(7, 7).Value = 93
$ws.Cells.Item(7, 8).Value = "◇サイト ○WordPress"

# Row 8
$ws.Cells.Item(8, 1).Value = "2025-11-26 12:38:38"
$ws.Cells.Item(8, 2).Value = "製造業向けMR業務支援アプリケーションの機能開発エンジニア募集(Unity/C#)"
$ws.Cells.Item(8, 3).Value = "システム開発"
$ws.Cells.Item(8, 4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(8, 5).Value = "期限情報なし"
$ws.Cells.Item(8, 6).Value = "https://www.lancers.jp/work/detail/5441557"
$ws.Cells.Item(8, 7).Value = 93
$ws.Cells.Item(8, 8).Value = "◆開発 ◇アプリ"

# Row 9
$ws.Cells.Item(9, 1).Value = "2025-11-26 12:38:38"
$ws.Cells.Item(9, 2).Value = "【急募】縫製工場向けPL・CF可視化アプリのMVP開発"
$ws.Cells.Item(9, 3).Value = "システム開発"
$ws.Cells.Item(9, 4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(9, 5).Value = "期限情報なし"
$ws.Cells.Item(9, 6).Value = "https://www.lancers.jp/work/detail/5440957"
$ws.Cells.Item(9, 7).Value = 93
$ws.Cells.Item(9, 8).Value = "◆開発 ◇アプリ"

# Row 10
$ws.Cells.Item(10, 1).Value = "2025-11-26 12:38:38"
$ws.Cells.Item(10, 2).Value = "製造業向け 技能習得・作業トレーニングVRシステムの開発(Unity/R3)"
$ws.Cells.Item(10, 3).Value = "システム開発"
$ws.Cells.Item(10, 4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(10, 5).Value = "期限情報なし"
$ws.Cells.Item(10, 6).Value = "https://www.lancers.jp/work/detail/5441568"
$ws.Cells.Item(10, 7).Value = 83
$ws.Cells.Item(10, 8).Value = "◆開発"

# Row 11
$ws.Cells.Item(11, 1).Value = "2025-11-26 12:38:38"
$ws.Cells.Item(11, 2).Value = "【急募】出会いサイトのテストユーザを募集します!"
$ws.Cells.Item(11, 3).Value = "システム開発"
$ws.Cells.Item(11, 4).Value = "1,000 ~ 5,000 円 / 固定"
$ws.Cells.Item(11, 5).Value = "期限情報なし"
$ws.Cells.Item(11, 6).Value = "https://www.lancers.jp/work/detail/5441448"
$ws.Cells.Item(11, 7).Value = 30
$ws.Cells.Item(11, 8).Value = "◇サイト"

# Row 12
$ws.Cells.Item(12, 1).Value = "2025-11-26 12:38:38"
$ws.Cells.Item(12, 2).Value = "【急募】出会いサイトのテストユーザを募集します!"
$ws.Cells.Item(12, 3).Value = "システム開発"
$ws.Cells.Item(12, 4).Value = "1,000 ~ 5,000 円 / 固定"
$ws.Cells.Item(12, 5).Value = "期限情報なし"
$ws.Cells.Item(12, 6).Value = "https://www.lancers.jp/work/detail/5441440"
$ws.Cells.Item(12, 7).Value = 30
$ws.Cells.Item(12, 8).Value = "◇サイト"

# Row 13
$ws.Cells.Item(13, 1).Value = "2025-11-26 12:38:38"
$ws.Cells.Item(13, 2).Value = "急募 限定公開 PR 限定公開の仕事"
$ws.Cells.Item(13, 3).Value = "システム開発"
$ws.Cells.Item(13, 4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(13, 5).Value = "期限情報なし"
$ws.Cells.Item(13, 6).Value = "https://www.lancers.jp/work/detail/5440230"
$ws.Cells.Item(13, 7).Value = 25

# Row 14
$ws.Cells.Item(14, 1).Value = "2025-11-26 12:38:38"
$ws.Cells.Item(14, 2).Value = "【急募】Salesforce設定・構築経験者を求む!高報酬案件"
$ws.Cells.Item(14, 3).Value = "システム開発"
$ws.Cells.Item(14, 4).Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Cells.Item(14, 5).Value = "期限情報なし"
$ws.Cells.Item(14, 6).Value = "https://www.lancers.jp/work/detail/5441609"
$ws.Cells.Item(14, 7).Value = 25

# Row 15
$ws.Cells.Item(15, 1).Value = "2025-11-26 12:38:38"
$ws.Cells.Item(15, 2).Value = "【急募】大規模プロジェクト統括のプロジェクトマネージャー募集(月:80万円~120万円)"
$ws.Cells.Item(15, 3).Value = "システム開発"
$ws.Cells.Item(15, 4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(15, 5).Value = "期限情報なし"
$ws.Cells.Item(15, 6).Value = "https://www.lancers.jp/work/detail/5441422"
$ws.Cells.Item(15, 7).Value = 25

# Re-create the hyperlinks on column F in row order (rId1..rId14) and make sure
# the cells keep the original "Hyperlink" cell style (same xf used before the edit)
# instead of the fresh style Hyperlinks.Add would normally stamp on.
$styleSource = $ws.Range("F2")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5441470")
$styleSource.Copy()
$ws.Range("F2").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5217096")
$styleSource.Copy()
$ws.Range("F3").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5441252")
$styleSource.Copy()
$ws.Range("F4").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5441207")
$styleSource.Copy()
$ws.Range("F5").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5441442")
$styleSource.Copy()
$ws.Range("F6").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5441612")
$styleSource.Copy()
$ws.Range("F7").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5441557")
$styleSource.Copy()
$ws.Range("F8").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5440957")
$styleSource.Copy()
$ws.Range("F9").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5441568")
$styleSource.Copy()
$ws.Range("F10").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5441448")
$styleSource.Copy()
$ws.Range("F11").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5441440")
$styleSource.Copy()
$ws.Range("F12").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("F13"), "https://www.lancers.jp/work/detail/5440230")
$styleSource.Copy()
$ws.Range("F13").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("F14"), "https://www.lancers.jp/work/detail/5441609")
$styleSource.Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("F15"), "https://www.lancers.jp/work/detail/5441422")
$styleSource.Copy()
$ws.Range("F15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column H got wider to fit the longer skill-summary strings now present
$ws.Columns.Item(8).ColumnWidth = 16.166666666666668